$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "27.792.41"
$ws.Range("E2").Value = "  -2.40%  "

# Row 3
$ws.Range("D3").Value = "1.782.21"
$ws.Range("E3").Value = "  -2.05%  "

# Row 4
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.002"
$ws.Range("E4").Value = "  -0.03%  "

# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "310.85"
$ws.Range("E5").Value = "  -1.89%  "

# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5116"
$ws.Range("E7").Value = "  -0.93%  "

# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3808"
$ws.Range("E8").Value = "  -2.05%  "

# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07778"
$ws.Range("E9").Value = "  -8.14%  "

# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "41.02"
$ws.Range("E10").Value = "  -1.86%  "

# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "1.086"
$ws.Range("E11").Value = "  -2.28%  "

# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.002"

# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "6.190"
$ws.Range("E13").Value = "  -3.81%  "

# Row 15
$ws.Range("D15").Value = "1.776.18"
$ws.Range("E15").Value = "  -2.32%  "

# Row 16
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "7.180"
$ws.Range("E16").Value = "  -4.42%  "

# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "91.28"
$ws.Range("E17").Value = "  -1.63%  "

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.00001074"
$ws.Range("E18").Value = "  -5.95%  "

# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06550"
$ws.Range("E19").Value = "  -1.41%  "

# Row 20
$ws.Range("E20").Value = "  +0.03%  "

# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "17.01"
$ws.Range("E21").Value = "  -4.20%  "

# Row 22
$ws.Range("E22").Value = "  -3.01%  "

# Row 23
$ws.Range("D23").Value = "27.833.33"
$ws.Range("E23").Value = "  -2.38%  "

# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.00"
$ws.Range("E24").Value = "  -3.91%  "

# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.233"
$ws.Range("E25").Value = "  -1.70%  "

# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "159.50"
$ws.Range("E26").Value = "  +0.14%  "

# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "20.18"
$ws.Range("E27").Value = "  -4.08%  "

# Row 28
$ws.Range("D28").Value = "1.983.38"
$ws.Range("E28").Value = "  -2.30%  "

# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.357"
$ws.Range("E29").Value = "  -1.70%  "

# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "123.69"
$ws.Range("E30").Value = "  -1.55%  "

# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1069"
$ws.Range("E31").Value = "  -1.84%  "

# Row 32
$ws.Range("E32").Value = "  -5.87%  "

# Row 33
$ws.Range("E33").Value = "  -0.44%  "

# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "5.479"
$ws.Range("E34").Value = "  -4.28%  "

# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.07080"
$ws.Range("E35").Value = "  -4.85%  "

# Row 36
$ws.Range("B36").Value = "FraxShare"
$ws.Range("C36").Value = "https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "8.816"
$ws.Range("E36").Value = "  -0.34%  "

# Row 37
$ws.Range("B37").Value = "VeChain"
$ws.Range("C37").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02304"
$ws.Range("E37").Value = "  -2.35%  "

# Row 38
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2117"
$ws.Range("E38").Value = "  -5.22%  "

# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "11.46"

# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "4.990"
$ws.Range("E40").Value = "  -4.15%  "

# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6085"
$ws.Range("E41").Value = "  -3.70%  "

# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.002"
$ws.Range("E42").Value = "  +0.09%  "

# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.151"
$ws.Range("E43").Value = "  -3.59%  "

# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.324"
$ws.Range("E44").Value = "  -5.51%  "

# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.04"
$ws.Range("E45").Value = "  -3.41%  "

# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.5885"
$ws.Range("E46").Value = "  -1.00%  "

# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.706"
$ws.Range("E47").Value = "  -2.03%  "

# Row 48
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "125.98"
$ws.Range("E48").Value = "  -0.22%  "

# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.192"
$ws.Range("E49").Value = "  -0.90%  "

# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.894"
$ws.Range("E50").Value = "  -4.73%  "

# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06841"
$ws.Range("E51").Value = "  -1.90%  "
